$p = $ppt.ActivePresentation

# 1. Remove the custom yellow fill override from the triangle shape on
#    slide 3, reverting it back to its shape-style (theme) fill. There is
#    no direct "clear override" COM property, so we pick up the (unfilled)
#    format of a freshly-added default shape and apply it to the triangle,
#    then discard the temporary donor shape.
$s3 = $p.Slides.Item(3)
$triangle = $null
for ($i = 1; $i -le $s3.Shapes.Count; $i++) {
    $shp = $s3.Shapes.Item($i)
    if ($shp.Name -eq "Isosceles Triangle 3") {
        $triangle = $shp
    }
}
$donor = $s3.Shapes.AddShape(7, 100, 100, 100, 100)
$donor.PickUp()
$triangle.Apply()
$donor.Delete()

# 2. Update the speaker notes on slide 3: drop the word "yellow".
$notes3 = $s3.NotesPage
for ($i = 1; $i -le $notes3.Shapes.Count; $i++) {
    $shp = $notes3.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        if ($shp.TextFrame.TextRange.Text -eq "Did you predict the next shape would be a yellow triangle?") {
            $shp.TextFrame.TextRange.Text = "Did you predict the next shape would be a triangle?"
        }
    }
}

# 3. Delete slide 4 (and its notes page) entirely - it was an empty,
#    unused slide.
$p.Slides.Item(4).Delete()
